$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referensi")

$ws.Range("B14").Value = "https://www.youtube.com/watch?v=xcn7hz7J7sI"
$ws.Range("A14").Value = "Jayanam"
$ws.Range("C14").Value = "Camera 3rd person mouse aim"
$ws.Range("D14").Value = "kamera arah tetikus"

$ws.Range("A15").Value = "Andrey Kubyshkin"
$ws.Range("B15").Value = "https://forum.unity.com/threads/moving-character-relative-to-camera.383086/"
$ws.Range("C15").Value = "Move character relative to camera"
$ws.Range("D15").Value = "Gerakan karakter relative dengan arah hadapan kamera"

$wsAlat = $wb.Worksheets.Item("Alat")
$wsAlat.Activate()
$wsAlat.Range("B12").Select() | Out-Null

$ws.Activate()
$ws.Range("D15").Select() | Out-Null
